$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TSLA")

# Row 8 - Interest Expense (Operating)
$ws.Range("B8").Value = 678000000.0

# Row 9 - Non-operating Income/Expense
$ws.Range("B9").Value = -187000000.0
$ws.Range("D9").Value = -269000000.0
$ws.Range("E9").Value = -255000000.0
$ws.Range("F9").Value = -449000000.0
$ws.Range("G9").Value = -596000000.0

# Row 10 - Non-operating Interest Expenses
$ws.Range("B10").Value = 30000000.0
$ws.Range("D10").Value = 194000000.0
$ws.Range("E10").Value = 373000000.0
$ws.Range("F10").Value = 537000000.0
$ws.Range("G10").Value = 685000000.0

# Row 15 - EPS (Basic)
$ws.Range("B15").Value = 1.0
$ws.Range("D15").Value = 0.506
$ws.Range("E15").Value = 0.396
$ws.Range("F15").Value = -0.166
$ws.Range("G15").Value = -1.006

# Row 21 - EBITDA
$ws.Range("B21").Value = 4594000000.0
$ws.Range("D21").Value = 4109000000.0
$ws.Range("E21").Value = 3563000000.0
$ws.Range("F21").Value = 3123000000.0
$ws.Range("G21").Value = 2273000000.0

# Row 23 - EPS (Diluted, from Cont. Ops)
$ws.Range("C23").Value = 0.7959
$ws.Range("D23").Value = 0.7023
$ws.Range("E23").Value = 0.5311
$ws.Range("F23").Value = -0.033

# Row 28 - EBITDA Margin
$ws.Range("B28").Value = 0.1278
